$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark rows 11 and 13 ("Done?" column G) as completed: "No" -> "Yes",
# using Excel's built-in "Good" cell style (green fill / green text),
# the same way the existing "No" entries use the built-in "Bad" style.
$ws.Range("G11").Value = "Yes"
$ws.Range("G11").Style = "Good"

$ws.Range("G13").Value = "Yes"
$ws.Range("G13").Style = "Good"

# Move the saved selection to D17, matching the updated cursor position.
$ws.Range("D17").Select()
